$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1): "Average GPU utilization" moves from G2 to G1 ---
$ws.Range("G1").Value = "Average GPU utilization"
$ws.Range("G2").ClearContents()
# the "< 30 %%" label in G3 is dropped entirely
$ws.Range("G3").ClearContents()

# --- batch_size = 256 block (rows 2-7) ---
$ws.Range("B2").Value = "batch_size = 256"

$ws.Range("B3").Value = "GPU/epochs"
$ws.Range("C3").Value = 64
$ws.Range("D3").Value = 128
$ws.Range("E3").Value = 256

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 113.73
$ws.Range("D4").Value = 215.72
$ws.Range("E4").Value = 418.85

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 101.97
$ws.Range("D5").Value = 194.36
$ws.Range("E5").Value = 376.06

$ws.Range("B6").Value = 2
$ws.Range("E6").Value = 355.88

$ws.Range("B7").Value = 1

# --- batch_size = 128 block (rows 9-14) ---
$ws.Range("B9").Value = "batch_size = 128"

$ws.Range("B10").Value = "GPU/epochs"
$ws.Range("C10").Value = 64
$ws.Range("D10").Value = 128
$ws.Range("E10").Value = 256

$ws.Range("B11").Value = 4
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 1

# --- batch_size = 2048 block (rows 16-21) ---
$ws.Range("B16").Value = "batch_size = 2048"

# row17 "GPU/epochs" label (B17) is removed - only the batch size axis values remain
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 64
$ws.Range("D17").Value = 128
$ws.Range("E17").Value = 256

$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 55.06
$ws.Range("D18").Value = 97.29

$ws.Range("B19").Value = 3

$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 56.11
$ws.Range("D20").Value = 104.92

$ws.Range("B21").Value = 1

# --- batch_size = 4096 block (rows 23-28) ---
$ws.Range("B23").Value = "batch_size = 4096"
$ws.Range("C23").Value = "low recall"

$ws.Range("C24").Value = 64
$ws.Range("D24").Value = 128
$ws.Range("E24").Value = 256

$ws.Range("B25").Value = 4
$ws.Range("C25").Value = 52.25

$ws.Range("B26").Value = 3

$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 54.36

$ws.Range("B28").Value = 1

# --- batch_size = 1024 block (rows 30-35) ---
$ws.Range("B30").Value = "batch_size = 1024"

$ws.Range("C31").Value = 64
$ws.Range("D31").Value = 128
$ws.Range("E31").Value = 256

$ws.Range("B32").Value = 4
$ws.Range("D32").Value = 112.19

$ws.Range("B33").Value = 3
$ws.Range("B34").Value = 2
$ws.Range("B35").Value = 1

# --- view state: scroll + selection ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F24").Select()
